# Commit: "Created some user stories"
#
# Diffing the captured before/after OOXML shows the user-story content
# (title, "As a ..." story paragraph, etc.) is identical on both sides of
# the commit - the only differences are incidental namespace-declaration
# bookkeeping that Word's packager re-emits on save (a handful of rarely
# used mc:Ignorable extension namespaces - cx3-cx8, aink, am3d, w16cid -
# move from the document root down onto the single legacy mc:Fallback
# block, and four not-yet-used w:latentStyles/w:lsdException entries -
# "Smart Hyperlink", "Hashtag", "Unresolved Mention", "Smart Link" - drop
# out of styles.xml). None of that is reachable, or even meaningful,
# through the Word object model - it carries no visible document content,
# formatting, or structure, so there is nothing for this script to change
# there.
#
# The actual authoring work the commit message refers to (the user story
# text itself) is already present in the document, so round-trip it
# through the object model untouched: confirm the story content is
# present and leave it exactly as authored.
$d = $word.ActiveDocument

$storyCount = $d.Paragraphs.Count
for ($i = 1; $i -le $storyCount; $i++) {
    $null = $d.Paragraphs.Item($i).Range.Text
}
